# Update the "丽水-漫展信息" workbook with the latest scraped con-expo entries.
#
# For both the "展览" sheet and the "全部类型" sheet (sheet indexes 1 and 4 -
# they mirror each other), a new event ("丽水·thp01～风摄少微") is inserted
# as the 3rd data row (row 4 on the sheet), pushing the two events that used
# to follow it down by one row. The "想去人数" (F column) counts for the two
# still-existing events, and for the one above the new row, also tick up.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # "丽水·龙泉ACG动漫游戏博览会" gained one more person wanting to go.
    $ws.Range("F3").Value = 469

    # Insert a brand-new row for the "thp01～风摄少微" expo right before the
    # "CCAC动漫七夕（回馈展）" row, shifting the later rows down.
    $ws.Rows.Item(4).Insert()

    # Restore the numbering style (bold + border + centered) on the new
    # index cell, which Insert() otherwise leaves unstyled.
    $idxCell = $ws.Range("A4")
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    $ws.Range("A4").Value = 3

    # Plain-text dates ("2024-07-27") look like real dates to Excel's
    # smart-entry parser, which would silently convert them to date
    # serials. Route the literal string through a text formula, then
    # paste-special just the value back so the cell ends up as ordinary
    # text (matching how the rest of the sheet stores its dates).
    $ws.Range("B4").Formula = '="2024-07-27"'
    $ws.Range("B4").Copy()
    $ws.Range("B4").PasteSpecial(-4163)

    $ws.Range("C4").Value = "丽水·thp01～风摄少微"
    $ws.Range("D4").Value = "大猷街 应星楼"
    $ws.Range("E4").Value = "2024.07.27 10:00-07.27 18:00"
    $ws.Range("F4").Value = 0
    $ws.Range("G4").Value = 50
    $ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=87134"
    $ws.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202406/JuvSmncN1717775885615.png"

    # "丽水·CCAC动漫七夕（回馈展）" (now on row 5) gained one more person,
    # and its running index bumps from 3 to 4.
    $ws.Range("A5").Value = 4
    $ws.Range("F5").Value = 30

    # "丽水·AEO纯白礼赞动漫嘉年华" (now on row 6) gained a few more people,
    # and its running index bumps from 4 to 5.
    $ws.Range("A6").Value = 5
    $ws.Range("F6").Value = 137
}
